# Update "CursosBB-2021" worksheet: mark a few courses as done ("V") in
# column C, move the active selection to C3, and drop the extra row height
# that had been set on row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "FEITO" column with "V" for the completed courses.
$ws.Range("C2").Value = "V"
$ws.Range("C3").Value = "V"
$ws.Range("C8").Value = "V"
$ws.Range("C10").Value = "V"

# Row 8 no longer needs the taller (30pt) row height; restore auto height.
$ws.Rows.Item(8).AutoFit()

# Update the active cell/selection.
$ws.Range("C3").Select()
